$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.683.96"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'  -3.06%  "
$ws.Range("E2").ClearFormats()
$ws.Range("D3").Value = "'2.621.44"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'  -1.17%  "
$ws.Range("E3").ClearFormats()
$ws.Range("E4").Value = "'  +0.07%  "
$ws.Range("E4").ClearFormats()
$ws.Range("D5").Value = "'574.45"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'  -3.96%  "
$ws.Range("E5").ClearFormats()
$ws.Range("D6").Value = "'155.14"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "'  -1.00%  "
$ws.Range("E6").ClearFormats()
$ws.Range("E7").Value = "'  +0.04%  "
$ws.Range("E7").ClearFormats()
$ws.Range("D8").Value = "'0.630"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "'  +0.49%  "
$ws.Range("E8").ClearFormats()
$ws.Range("E9").Value = "'  -5.09%  "
$ws.Range("E9").ClearFormats()
$ws.Range("E10").Value = "'  -0.41%  "
$ws.Range("E10").ClearFormats()
$ws.Range("E11").Value = "'  -3.07%  "
$ws.Range("E11").ClearFormats()
$ws.Range("E12").Value = "'  -0.38%  "
$ws.Range("E12").ClearFormats()
$ws.Range("D13").Value = "'28.28"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "'  -1.37%  "
$ws.Range("E13").ClearFormats()
$ws.Range("D14").Value = "'3.095.82"
$ws.Range("D14").ClearFormats()
$ws.Range("D15").Value = "'0.0000184"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "'  -7.00%  "
$ws.Range("E15").ClearFormats()
$ws.Range("D16").Value = "'63.632.57"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "'  -2.86%  "
$ws.Range("E16").ClearFormats()
$ws.Range("D17").Value = "'2.607.36"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "'  -2.48%  "
$ws.Range("E17").ClearFormats()
$ws.Range("D18").Value = "'12.07"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "'  -4.29%  "
$ws.Range("E18").ClearFormats()
$ws.Range("E19").Value = "'  -2.00%  "
$ws.Range("E19").ClearFormats()
$ws.Range("D20").Value = "'7.52"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "'  +1.10%  "
$ws.Range("E20").ClearFormats()
$ws.Range("D21").Value = "'343.83"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "'  -1.42%  "
$ws.Range("E21").ClearFormats()
$ws.Range("E22").Value = "'  -0.03%  "
$ws.Range("E22").ClearFormats()
$ws.Range("D23").Value = "'67.07"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "'  -2.74%  "
$ws.Range("E23").ClearFormats()
$ws.Range("E24").Value = "'  +3.20%  "
$ws.Range("E24").ClearFormats()
$ws.Range("E25").Value = "'  -3.88%  "
$ws.Range("E25").ClearFormats()
$ws.Range("D26").Value = "'9.24"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "'  -4.23%  "
$ws.Range("E26").ClearFormats()
$ws.Range("D27").Value = "'583.82"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "'  +9.65%  "
$ws.Range("E27").ClearFormats()
$ws.Range("E28").Value = "'  +0.20%  "
$ws.Range("E28").ClearFormats()
$ws.Range("B29").Value = "'Binance-PegBSC-USD"
$ws.Range("B29").ClearFormats()
$ws.Range("C29").Value = "'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("C29").ClearFormats()
$ws.Range("D29").Value = "'1.00"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "'  +0.20%  "
$ws.Range("E29").ClearFormats()
$ws.Range("B30").Value = "'Kaspa"
$ws.Range("B30").ClearFormats()
$ws.Range("C30").Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("C30").ClearFormats()
$ws.Range("D30").Value = "'0.161"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "'  -2.18%  "
$ws.Range("E30").ClearFormats()
$ws.Range("E31").Value = "'  -0.16%  "
$ws.Range("E31").ClearFormats()
$ws.Range("D32").Value = "'2.07"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "'  -2.65%  "
$ws.Range("E32").ClearFormats()
$ws.Range("D33").Value = "'1.70"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "'  -3.56%  "
$ws.Range("E33").ClearFormats()
$ws.Range("D34").Value = "'6.52"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "'  +1.46%  "
$ws.Range("E34").ClearFormats()
$ws.Range("D35").Value = "'5.32"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "'  -1.68%  "
$ws.Range("E35").ClearFormats()
$ws.Range("D36").Value = "'0.409"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "'  -2.31%  "
$ws.Range("E36").ClearFormats()
$ws.Range("D37").Value = "'19.88"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "'  -2.43%  "
$ws.Range("E37").ClearFormats()
$ws.Range("E38").Value = "'  +0.10%  "
$ws.Range("E38").ClearFormats()
$ws.Range("D39").Value = "'153.46"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "'  -1.56%  "
$ws.Range("E39").ClearFormats()
$ws.Range("E40").Value = "'  -3.39%  "
$ws.Range("E40").ClearFormats()
$ws.Range("E41").Value = "'  -0.06%  "
$ws.Range("E41").ClearFormats()
$ws.Range("D42").Value = "'41.33"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "'  -2.64%  "
$ws.Range("E42").ClearFormats()
$ws.Range("D43").Value = "'156.82"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "'  -2.82%  "
$ws.Range("E43").ClearFormats()
$ws.Range("D44").Value = "'2.38"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "'  +4.20%  "
$ws.Range("E44").ClearFormats()
$ws.Range("D45").Value = "'3.95"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "'  -2.69%  "
$ws.Range("E45").ClearFormats()
$ws.Range("E46").Value = "'  -1.96%  "
$ws.Range("E46").ClearFormats()
$ws.Range("D47").Value = "'22.82"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "'  +0.91%  "
$ws.Range("E47").ClearFormats()
$ws.Range("D48").Value = "'0.631"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "'  -0.64%  "
$ws.Range("E48").ClearFormats()
$ws.Range("D49").Value = "'0.101"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "'  +1.87%  "
$ws.Range("E49").ClearFormats()
$ws.Range("D50").Value = "'0.0251"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "'  -1.12%  "
$ws.Range("E50").ClearFormats()
$ws.Range("D51").Value = "'19.06"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "'  -3.46%  "
$ws.Range("E51").ClearFormats()
